# chore: update Sheets via scheduled runner
# Refresh cached market-board price / profit figures (columns H-N) for
# affected leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3237.4375
$ws.Range("I100").Value = 2600
$ws.Range("J100").Value = 3874.875
$ws.Range("K100").Value = 2600
$ws.Range("L100").Value = 3874.875
$ws.Range("M100").Value = -2059
$ws.Range("N100").Value = -4956.875

$ws.Range("H137").Value = 5558203
$ws.Range("I137").Value = 1650.591
$ws.Range("K137").Value = 4951.772999999999
$ws.Range("M137").Value = -2401.772999999999

$ws.Range("H138").Value = 4632268
$ws.Range("I138").Value = 1459.72
$ws.Range("J138").Value = 8624344
$ws.Range("K138").Value = 4379.16
$ws.Range("L138").Value = 25873032
$ws.Range("M138").Value = 760.8400000000001
$ws.Range("N138").Value = -25883312

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1930.1428
$ws.Range("I2").Value = 1477.75
$ws.Range("J2").Value = 2533.3333
$ws.Range("K2").Value = 1477.75
$ws.Range("L2").Value = 2533.3333
$ws.Range("M2").Value = -1364.75
$ws.Range("N2").Value = -2759.3333

$ws.Range("H45").Value = 3200.375
$ws.Range("I45").Value = 3267.1667
$ws.Range("K45").Value = 3267.1667
$ws.Range("M45").Value = -2890.1667

$ws.Range("H88").Value = 2724.625
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2724.625
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 2724.625
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -3536.625

$ws.Range("H91").Value = 2724.625
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2724.625
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 2724.625
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -5532.625

$ws.Range("H110").Value = 1400
$ws.Range("I110").Value = 1400
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1400
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 645
$ws.Range("N110").ClearContents()

$ws.Range("H116").Value = 1930.1428
$ws.Range("I116").Value = 1477.75
$ws.Range("J116").Value = 2533.3333
$ws.Range("K116").Value = 1477.75
$ws.Range("L116").Value = 2533.3333
$ws.Range("M116").Value = 816.25
$ws.Range("N116").Value = -7121.3333

$ws.Range("H134").Value = 49060.715
$ws.Range("J134").Value = 49060.715
$ws.Range("L134").Value = 49060.715
$ws.Range("N134").Value = -59200.715

$ws.Range("H135").Value = 19993.375
$ws.Range("J135").Value = 19993.375
$ws.Range("L135").Value = 19993.375
$ws.Range("N135").Value = -30133.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1930.1428
$ws.Range("I3").Value = 1477.75
$ws.Range("J3").Value = 2533.3333
$ws.Range("K3").Value = 1477.75
$ws.Range("L3").Value = 2533.3333
$ws.Range("M3").Value = -1363.75
$ws.Range("N3").Value = -2761.3333

$ws.Range("H86").Value = 22729156
$ws.Range("I86").Value = 1904.7142
$ws.Range("J86").Value = 62501850
$ws.Range("K86").Value = 1904.7142
$ws.Range("L86").Value = 62501850
$ws.Range("M86").Value = -781.7141999999999
$ws.Range("N86").Value = -62504096

$ws.Range("H89").Value = 22729156
$ws.Range("I89").Value = 1904.7142
$ws.Range("J89").Value = 62501850
$ws.Range("K89").Value = 9523.571
$ws.Range("L89").Value = 312509250
$ws.Range("M89").Value = -3907.571
$ws.Range("N89").Value = -312520482

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11499912
$ws.Range("I31").Value = 13701.7
$ws.Range("J31").Value = 17545286
$ws.Range("K31").Value = 13701.7
$ws.Range("L31").Value = 17545286
$ws.Range("M31").Value = -13406.7
$ws.Range("N31").Value = -17545876

$ws.Range("H34").Value = 11499912
$ws.Range("I34").Value = 13701.7
$ws.Range("J34").Value = 17545286
$ws.Range("K34").Value = 13701.7
$ws.Range("L34").Value = 17545286
$ws.Range("M34").Value = -13499.7
$ws.Range("N34").Value = -17545690

$ws.Range("H62").Value = 2364.8
$ws.Range("I62").Value = 2288.3333
$ws.Range("J62").Value = 2670.6667
$ws.Range("K62").Value = 2288.3333
$ws.Range("L62").Value = 2670.6667
$ws.Range("M62").Value = -1664.3333
$ws.Range("N62").Value = -3918.6667

$ws.Range("H65").Value = 2364.8
$ws.Range("I65").Value = 2288.3333
$ws.Range("J65").Value = 2670.6667
$ws.Range("K65").Value = 11441.6665
$ws.Range("L65").Value = 13353.3335
$ws.Range("M65").Value = -8321.666499999999
$ws.Range("N65").Value = -19593.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 732.4032
$ws.Range("I113").Value = 601.6042
$ws.Range("J113").Value = 1180.8572
$ws.Range("K113").Value = 1804.8126
$ws.Range("L113").Value = 3542.5716
$ws.Range("M113").Value = 365.1874
$ws.Range("N113").Value = -7882.571599999999

$ws.Range("H118").Value = 1607.32
$ws.Range("J118").Value = 1801.1904
$ws.Range("L118").Value = 5403.5712
$ws.Range("N118").Value = -7889.5712

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2282.8
$ws.Range("I113").Value = 1800.3334
$ws.Range("K113").Value = 1800.3334
$ws.Range("M113").Value = 369.6666

$ws.Range("H122").Value = 5131828
$ws.Range("I122").Value = 11114478
$ws.Range("J122").Value = 3842.8572
$ws.Range("K122").Value = 33343434
$ws.Range("L122").Value = 11528.5716
$ws.Range("M122").Value = -33340984
$ws.Range("N122").Value = -16428.5716

$ws.Range("H126").Value = 4668.25
$ws.Range("I126").Value = 3095.4
$ws.Range("J126").Value = 5192.533
$ws.Range("K126").Value = 9286.200000000001
$ws.Range("L126").Value = 15577.599
$ws.Range("M126").Value = -6816.200000000001
$ws.Range("N126").Value = -20517.599

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4585.3076
$ws.Range("I16").Value = 3300.75
$ws.Range("J16").Value = 20000
$ws.Range("K16").Value = 3300.75
$ws.Range("L16").Value = 20000
$ws.Range("M16").Value = -3130.75
$ws.Range("N16").Value = -20340

$ws.Range("H46").Value = 865.55554
$ws.Range("I46").Value = 727.1429000000001
$ws.Range("J46").Value = 1350
$ws.Range("K46").Value = 727.1429000000001
$ws.Range("L46").Value = 1350
$ws.Range("M46").Value = -539.1429000000001
$ws.Range("N46").Value = -1726

$ws.Range("H68").Value = 2080
$ws.Range("I68").Value = 1500
$ws.Range("J68").Value = 2196
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 2196
$ws.Range("M68").Value = -751
$ws.Range("N68").Value = -3694

$ws.Range("H71").Value = 2080
$ws.Range("I71").Value = 1500
$ws.Range("J71").Value = 2196
$ws.Range("K71").Value = 7500
$ws.Range("L71").Value = 10980
$ws.Range("M71").Value = -3756
$ws.Range("N71").Value = -18468
